$d = $word.ActiveDocument
$d.Content.Find.Execute("2023-05-21 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-05-22 Monday", 2)
$d.Content.Find.Execute("62-5=", $true, $false, $false, $false, $false, $true, 1, $false, "74-6=", 2)
$d.Content.Find.Execute("66-56=", $true, $false, $false, $false, $false, $true, 1, $false, "2+76=", 2)
$d.Content.Find.Execute("27-9=", $true, $false, $false, $false, $false, $true, 1, $false, "84-78=", 2)
$d.Content.Find.Execute("91-26=", $true, $false, $false, $false, $false, $true, 1, $false, "66+8=", 2)
$d.Content.Find.Execute("81-21=", $true, $false, $false, $false, $false, $true, 1, $false, "99-99=", 2)
$d.Content.Find.Execute("29-9=", $true, $false, $false, $false, $false, $true, 1, $false, "44-32=", 2)
$d.Content.Find.Execute("13+86=", $true, $false, $false, $false, $false, $true, 1, $false, "67-20=", 2)
$d.Content.Find.Execute("58-46=", $true, $false, $false, $false, $false, $true, 1, $false, "38-0=", 2)
$d.Content.Find.Execute("54-42=", $true, $false, $false, $false, $false, $true, 1, $false, "35+3=", 2)
$d.Content.Find.Execute("27+41=", $true, $false, $false, $false, $false, $true, 1, $false, "68-27=", 2)
$d.Content.Find.Execute("1+66=", $true, $false, $false, $false, $false, $true, 1, $false, "36+20=", 2)
$d.Content.Find.Execute("75-33=", $true, $false, $false, $false, $false, $true, 1, $false, "95-35=", 2)
$d.Content.Find.Execute("21+32=", $true, $false, $false, $false, $false, $true, 1, $false, "6+0=", 2)
$d.Content.Find.Execute("84-7=", $true, $false, $false, $false, $false, $true, 1, $false, "70-38=", 2)
$d.Content.Find.Execute("47+26=", $true, $false, $false, $false, $false, $true, 1, $false, "48-7=", 2)
$d.Content.Find.Execute("52-11=", $true, $false, $false, $false, $false, $true, 1, $false, "49-4=", 2)
$d.Content.Find.Execute("35-6=", $true, $false, $false, $false, $false, $true, 1, $false, "97-69=", 2)
$d.Content.Find.Execute("94-76=", $true, $false, $false, $false, $false, $true, 1, $false, "35-8=", 2)
$d.Content.Find.Execute("64+6=", $true, $false, $false, $false, $false, $true, 1, $false, "31-3=", 2)
$d.Content.Find.Execute("72-67=", $true, $false, $false, $false, $false, $true, 1, $false, "68-21=", 2)
$d.Content.Find.Execute("60+11=", $true, $false, $false, $false, $false, $true, 1, $false, "92-1=", 2)
$d.Content.Find.Execute("0+87=", $true, $false, $false, $false, $false, $true, 1, $false, "2+69=", 2)
$d.Content.Find.Execute("79-59=", $true, $false, $false, $false, $false, $true, 1, $false, "5+4=", 2)
$d.Content.Find.Execute("74-2=", $true, $false, $false, $false, $false, $true, 1, $false, "65-37=", 2)
$d.Content.Find.Execute("58-28=", $true, $false, $false, $false, $false, $true, 1, $false, "68-26=", 2)
$d.Content.Find.Execute("60-32=", $true, $false, $false, $false, $false, $true, 1, $false, "23+5=", 2)
$d.Content.Find.Execute("41-32=", $true, $false, $false, $false, $false, $true, 1, $false, "55+8=", 2)
$d.Content.Find.Execute("11+37=", $true, $false, $false, $false, $false, $true, 1, $false, "40+38=", 2)
$d.Content.Find.Execute("10+59=", $true, $false, $false, $false, $false, $true, 1, $false, "20+62=", 2)
$d.Content.Find.Execute("20+68=", $true, $false, $false, $false, $false, $true, 1, $false, "39-20=", 2)
$d.Content.Find.Execute("78-32=", $true, $false, $false, $false, $false, $true, 1, $false, "97-15=", 2)
$d.Content.Find.Execute("47-36=", $true, $false, $false, $false, $false, $true, 1, $false, "13+22=", 2)
$d.Content.Find.Execute("21+72=", $true, $false, $false, $false, $false, $true, 1, $false, "69-19=", 2)
$d.Content.Find.Execute("0+64=", $true, $false, $false, $false, $false, $true, 1, $false, "83-38=", 2)
$d.Content.Find.Execute("71-58=", $true, $false, $false, $false, $false, $true, 1, $false, "69-14=", 2)
$d.Content.Find.Execute("45+15=", $true, $false, $false, $false, $false, $true, 1, $false, "86-67=", 2)
$d.Content.Find.Execute("21+9=", $true, $false, $false, $false, $false, $true, 1, $false, "3+85=", 2)
$d.Content.Find.Execute("52+28=", $true, $false, $false, $false, $false, $true, 1, $false, "14+73=", 2)
$d.Content.Find.Execute("28+31=", $true, $false, $false, $false, $false, $true, 1, $false, "19+73=", 2)
$d.Content.Find.Execute("35+44=", $true, $false, $false, $false, $false, $true, 1, $false, "28+62=", 2)
$d.Content.Find.Execute("12+55=", $true, $false, $false, $false, $false, $true, 1, $false, "22+61=", 2)
$d.Content.Find.Execute("4+93=", $true, $false, $false, $false, $false, $true, 1, $false, "39+39=", 2)
$d.Content.Find.Execute("2+57=", $true, $false, $false, $false, $false, $true, 1, $false, "61-49=", 2)
$d.Content.Find.Execute("83-45=", $true, $false, $false, $false, $false, $true, 1, $false, "5+67=", 2)
$d.Content.Find.Execute("18+80=", $true, $false, $false, $false, $false, $true, 1, $false, "15-6=", 2)
$d.Content.Find.Execute("84-53=", $true, $false, $false, $false, $false, $true, 1, $false, "54+28=", 2)
$d.Content.Find.Execute("85+6=", $true, $false, $false, $false, $false, $true, 1, $false, "6+67=", 2)
$d.Content.Find.Execute("33-33=", $true, $false, $false, $false, $false, $true, 1, $false, "63+14=", 2)
$d.Content.Find.Execute("36+16=", $true, $false, $false, $false, $false, $true, 1, $false, "87-5=", 2)
$d.Content.Find.Execute("55+37=", $true, $false, $false, $false, $false, $true, 1, $false, "14+69=", 2)
$d.Content.Find.Execute("5+10=", $true, $false, $false, $false, $false, $true, 1, $false, "12-10=", 2)
$d.Content.Find.Execute("45-22=", $true, $false, $false, $false, $false, $true, 1, $false, "93-77=", 2)
$d.Content.Find.Execute("30-0=", $true, $false, $false, $false, $false, $true, 1, $false, "30-4=", 2)
$d.Content.Find.Execute("36+50=", $true, $false, $false, $false, $false, $true, 1, $false, "53+33=", 2)
$d.Content.Find.Execute("71-15=", $true, $false, $false, $false, $false, $true, 1, $false, "45-28=", 2)
$d.Content.Find.Execute("62-20=", $true, $false, $false, $false, $false, $true, 1, $false, "90+6=", 2)
$d.Content.Find.Execute("30-15=", $true, $false, $false, $false, $false, $true, 1, $false, "41+14=", 2)
$d.Content.Find.Execute("37+4=", $true, $false, $false, $false, $false, $true, 1, $false, "67-58=", 2)
$d.Content.Find.Execute("83-9=", $true, $false, $false, $false, $false, $true, 1, $false, "96-45=", 2)
$d.Content.Find.Execute("71-53=", $true, $false, $false, $false, $false, $true, 1, $false, "38+15=", 2)
$d.Content.Find.Execute("79-39=", $true, $false, $false, $false, $false, $true, 1, $false, "46+42=", 2)
$d.Content.Find.Execute("60-41=", $true, $false, $false, $false, $false, $true, 1, $false, "29+48=", 2)
$d.Content.Find.Execute("52+35=", $true, $false, $false, $false, $false, $true, 1, $false, "97-50=", 2)
$d.Content.Find.Execute("26-23=", $true, $false, $false, $false, $false, $true, 1, $false, "9+61=", 2)
$d.Content.Find.Execute("51-42=", $true, $false, $false, $false, $false, $true, 1, $false, "37+24=", 2)
$d.Content.Find.Execute("1+45=", $true, $false, $false, $false, $false, $true, 1, $false, "63-35=", 2)
$d.Content.Find.Execute("71-39=", $true, $false, $false, $false, $false, $true, 1, $false, "27+64=", 2)
$d.Content.Find.Execute("43+42=", $true, $false, $false, $false, $false, $true, 1, $false, "50-7=", 2)
$d.Content.Find.Execute("5+26=", $true, $false, $false, $false, $false, $true, 1, $false, "61+30=", 2)
$d.Content.Find.Execute("87-25=", $true, $false, $false, $false, $false, $true, 1, $false, "16-0=", 2)
$d.Content.Find.Execute("72+4=", $true, $false, $false, $false, $false, $true, 1, $false, "1+69=", 2)
$d.Content.Find.Execute("31+40=", $true, $false, $false, $false, $false, $true, 1, $false, "61-32=", 2)
$d.Content.Find.Execute("32+53=", $true, $false, $false, $false, $false, $true, 1, $false, "74+2=", 2)
$d.Content.Find.Execute("97-83=", $true, $false, $false, $false, $false, $true, 1, $false, "86+13=", 2)
$d.Content.Find.Execute("15-9=", $true, $false, $false, $false, $false, $true, 1, $false, "41+2=", 2)
$d.Content.Find.Execute("27-1=", $true, $false, $false, $false, $false, $true, 1, $false, "3-2=", 2)
$d.Content.Find.Execute("17+6=", $true, $false, $false, $false, $false, $true, 1, $false, "85+8=", 2)
$d.Content.Find.Execute("27-3=", $true, $false, $false, $false, $false, $true, 1, $false, "93-8=", 2)
$d.Content.Find.Execute("66-22=", $true, $false, $false, $false, $false, $true, 1, $false, "97-54=", 2)
$d.Content.Find.Execute("94-48=", $true, $false, $false, $false, $false, $true, 1, $false, "50-39=", 2)
$d.Content.Find.Execute("2+10=", $true, $false, $false, $false, $false, $true, 1, $false, "73+0=", 2)
$d.Content.Find.Execute("24+8=", $true, $false, $false, $false, $false, $true, 1, $false, "55+8=", 2)
$d.Content.Find.Execute("8+88=", $true, $false, $false, $false, $false, $true, 1, $false, "75+23=", 2)
$d.Content.Find.Execute("56-9=", $true, $false, $false, $false, $false, $true, 1, $false, "18+7=", 2)
$d.Content.Find.Execute("21+8=", $true, $false, $false, $false, $false, $true, 1, $false, "22+65=", 2)
$d.Content.Find.Execute("54+3=", $true, $false, $false, $false, $false, $true, 1, $false, "2+9=", 2)
$d.Content.Find.Execute("30+41=", $true, $false, $false, $false, $false, $true, 1, $false, "66+9=", 2)
$d.Content.Find.Execute("60-39=", $true, $false, $false, $false, $false, $true, 1, $false, "47-9=", 2)
$d.Content.Find.Execute("92-10=", $true, $false, $false, $false, $false, $true, 1, $false, "75+23=", 2)
$d.Content.Find.Execute("47-4=", $true, $false, $false, $false, $false, $true, 1, $false, "89-71=", 2)
$d.Content.Find.Execute("67-61=", $true, $false, $false, $false, $false, $true, 1, $false, "21+70=", 2)
$d.Content.Find.Execute("21+55=", $true, $false, $false, $false, $false, $true, 1, $false, "23+15=", 2)
$d.Content.Find.Execute("24+10=", $true, $false, $false, $false, $false, $true, 1, $false, "41+16=", 2)
$d.Content.Find.Execute("15+65=", $true, $false, $false, $false, $false, $true, 1, $false, "30-25=", 2)
$d.Content.Find.Execute("22+17=", $true, $false, $false, $false, $false, $true, 1, $false, "35-30=", 2)
$d.Content.Find.Execute("5+46=", $true, $false, $false, $false, $false, $true, 1, $false, "44+13=", 2)
$d.Content.Find.Execute("68-28=", $true, $false, $false, $false, $false, $true, 1, $false, "86-17=", 2)
$d.Content.Find.Execute("1+54=", $true, $false, $false, $false, $false, $true, 1, $false, "13+52=", 2)
$d.Content.Find.Execute("21+67=", $true, $false, $false, $false, $false, $true, 1, $false, "4+43=", 2)
$d.Content.Find.Execute("54-40=", $true, $false, $false, $false, $false, $true, 1, $false, "85-9=", 2)
